$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.901.62'
$ws.Range('E2').Value = '  -3.09%  '
$ws.Range('D3').Value = '3.416.06'
$ws.Range('E3').Value = '  -5.37%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '183.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -10.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '535.04'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.33%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '3.413.78'
$ws.Range('E8').Value = '  -5.35%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -6.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.16'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.93%  '
$ws.Range('E12').Value = '  -10.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000258'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -10.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.43'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -6.04%  '
$ws.Range('D15').Value = '3.954.87'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').Value = '3.406.91'
$ws.Range('E17').Value = '  -5.46%  '
$ws.Range('D18').Value = '65.620.51'
$ws.Range('E18').Value = '  -3.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.69'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -6.67%  '
$ws.Range('E20').Value = '  -7.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.991'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -8.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '381.91'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.78'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -9.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.98'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -15.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.69'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.65%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.72'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.20%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.77'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.63'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '700.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.99'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.84'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -19.39%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.16'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.80%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.32'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -7.29%  '
$ws.Range('E35').Value = '  -6.13%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '37.13'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -12.11%  '
$ws.Range('E38').Value = '  -7.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E40').Value = '  -5.49%  '
$ws.Range('D41').Value = '2.896.83'
$ws.Range('E41').Value = '  -10.47%  '
$ws.Range('E42').Value = '  -13.24%  '
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('E44').Value = '  -5.01%  '
$ws.Range('D45').Value = '0.0₃0631'
$ws.Range('E45').Value = '  -18.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.39'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -14.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.127'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '135.16'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.89'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.61'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.34'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -23.75%  '
